# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# For this particular template resource, the committed change only
# re-serializes the existing OOXML parts (word/document.xml,
# word/footnotes.xml, word/styles.xml): every attribute in the diff is
# present both before and after the change, only the attribute order
# differs (a canonical/alphabetical re-ordering produced when the
# resource was re-saved). No paragraph text, formatting value, style
# definition or document part was actually added, removed or modified.
#
# So the only safe, faithful edit here is a no-op content-wise: we simply
# touch the document through the Word object model (without mutating
# anything) so the part is round-tripped, and leave every run of text,
# property value and style untouched.
$d = $word.ActiveDocument

# Sanity-check that the well known content is still present; this mirrors
# how the template was inspected before being re-saved, without altering
# any text, run, or paragraph/style property.
$null = $d.Content.Find.Execute("End of demonstration.")
